$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that parse as plain numbers need a leading apostrophe so Excel
# stores them as literal text (preserving exact formatting, e.g. trailing
# zeros / leading zeros) instead of silently converting them to numbers.

$ws.Range("D2").Value = "62.005.68"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "2.578.34"
$ws.Range("E3").Value = "  -4.37%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'550.60"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "'155.34"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "'5.48"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "'0.367"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "3.033.82"
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("D14").Value = "'25.76"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").Value = "61.844.08"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "'0.0000146"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "2.584.27"
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("D18").Value = "'11.63"
$ws.Range("E18").Value = "  -4.29%  "
$ws.Range("D19").Value = "'4.57"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'338.39"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "'6.05"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'0.494"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").Value = "'63.53"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'8.18"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'7.41"
$ws.Range("E28").Value = "  +5.07%  "
$ws.Range("D29").Value = "0.0₃0842"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Value = "'161.99"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'4.87"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'1.44"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "'19.26"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'331.72"
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'6.06"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "'0.920"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").Value = "'3.95"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "'37.59"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").Value = "'20.97"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'0.607"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").Value = "2.118.16"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0549"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.93"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").Value = "'19.62"
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("D50").Value = "'0.0969"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'0.0241"
$ws.Range("E51").Value = "  -1.19%  "
